$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.666286706924438
$ws.Range("B1").Value = 2.425266981124878
$ws.Range("C1").Value = 3.477336645126343
$ws.Range("D1").Value = 1.288235306739807
$ws.Range("E1").Value = 0.8190047144889832
